# Update the "想去人数" (interested count) figures on the "展览" and
# "全部类型" worksheets: row 2 (丽水·第四届HP国风动漫游戏嘉年华) 545 -> 546,
# and row 7 (丽水·AEO纯白礼赞动漫嘉年华) 780 -> 781.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 546
    $ws.Range("F7").Value = 781
}
